$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.6544228876835514, 1.2261136499182372, -0.20512550369447885, 2.661927039588663)"
$ws.Range("C2").Value = "NIG(1.0124955062342962, 0.8202285279714585, 3.648492238598287, 5.735027511786299)"
$ws.Range("D2").Value = "NCT(3.2171879923591855, 1.475521859453918, -0.7749043550901462, 2.3496400526336663)"
$ws.Range("E2").Value = "NIG(1.2105011444355092, 0.8165292680338047, 5.46267382296954, 5.7592349157177924)"
